$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values ---
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Shared formatting for B1 (bold font, thin box border, centered/top aligned) ---
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Borders.Weight = 2
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").HorizontalAlignment = -4108

# Reuse the exact same style for A2 by copying formats (avoids generating a
# duplicate/unused style entry in the workbook's style table).
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
